# ---------------------------------------------------------------------------
# Locate the target paragraph:
#   "In the view folder... create a subfolder called layouts and in this
#    layout folder move the layout.hbs in there "
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.StartsWith("In the view folder") -and $t -like "*layout.hbs in there*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate target paragraph"
}

$ellipsis = [char]0x2026
$enDash   = [char]0x2013

# ---------------------------------------------------------------------------
# Rewrite the paragraph's text without the trailing space. Deleting the
# existing text first (rather than just assigning .Text) avoids the old
# run keeping a now-unnecessary xml:space="preserve" attribute.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item($targetIndex)
$r = $p.Range
$null = $r.MoveEnd(1, -1)  # exclude the paragraph mark
$null = $r.Delete()
$r.InsertAfter("In the view folder" + $ellipsis + " create a subfolder called layouts and in this layout folder move the layout.hbs in there")

# ---------------------------------------------------------------------------
# Insert six new, still-empty paragraphs right after it first. Doing all of
# the InsertParagraphAfter calls before any text/formatting is applied
# avoids the newly-typed formatting state "leaking" forward into paragraphs
# created later.
# ---------------------------------------------------------------------------
$cur = $d.Paragraphs.Item($targetIndex).Range
for ($k = 0; $k -lt 6; $k++) {
    $cur.InsertParagraphAfter()
    $cur = $d.Paragraphs.Item($targetIndex + $k + 1).Range
}

$idxLayoutHbs   = $targetIndex + 1
$idxBelowBody   = $targetIndex + 2
$idxPartialLine = $targetIndex + 3
$idxComment     = $targetIndex + 4
$idxSaveRefresh = $targetIndex + 5
$idxTrailSpace  = $targetIndex + 6

$d.Paragraphs.Item($idxLayoutHbs).Range.Text   = "layout.hbs: "
$d.Paragraphs.Item($idxBelowBody).Range.Text   = "below the <body> element and above the {{{body}}} hook type in"
$d.Paragraphs.Item($idxComment).Range.Text     = "// specifying the partials I was to include " + $enDash + " in this case the header"
$d.Paragraphs.Item($idxSaveRefresh).Range.Text = "Save and refresh the page" + $ellipsis + " now we can see the html is being populated by bootstrap and header.hbs"
$d.Paragraphs.Item($idxTrailSpace).Range.Text  = " "

# ---------------------------------------------------------------------------
# Fill the "{{> header }}" paragraph last (it contains styled runs, and
# nothing should be inserted after it once it has been styled).
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs.Item($idxPartialLine)
$rng = $p4.Range
$rng.Collapse(1)
$paraStart = $rng.Start

$segments = @("{{>", " ", "header", " ", "}}", " ")
$pos = $paraStart
$starts = @()
foreach ($seg in $segments) {
    $starts += $pos
    $insertRng = $d.Range($pos, $pos)
    $insertRng.InsertAfter($seg)
    $pos = $pos + $seg.Length
}

# Style the first five runs ("{{>", " ", "header", " ", "}}") as
# TerminalCodeOutputChar; leave the final trailing space run unstyled.
for ($i = 0; $i -lt 5; $i++) {
    $segStart = $starts[$i]
    $segEnd = $segStart + $segments[$i].Length
    $sRng = $d.Range($segStart, $segEnd)
    $sRng.Style = "Terminal Code Output Char"
}

Write-Output "done"
